$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="67.396.16"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Value = "  -3.63%  "
$ws.Range("D3").Formula = '="3.706.50"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Value = "  -4.33%  "
$ws.Range("D5").Formula = '="596.04"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -2.50%  "
$ws.Range("D6").Formula = '="165.63"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -5.42%  "
$ws.Range("D7").Formula = '="3.706.98"'
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  -4.36%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Formula = '="0.528"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("D10").Formula = '="0.161"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  -3.26%  "
$ws.Range("D11").Formula = '="6.17"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  -4.72%  "
$ws.Range("E12").Value = "  -4.18%  "
$ws.Range("D13").Formula = '="37.62"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  -6.26%  "
$ws.Range("D14").Formula = '="0.0000241"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  -5.47%  "
$ws.Range("D15").Formula = '="4.329.37"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  -4.11%  "
$ws.Range("D16").Formula = '="3.707.73"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  -4.04%  "
$ws.Range("D17").Formula = '="67.516.90"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  -3.48%  "
$ws.Range("D18").Formula = '="17.50"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  +4.90%  "
$ws.Range("D19").Formula = '="7.17"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  -3.90%  "
$ws.Range("D20").Formula = '="0.115"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  -2.97%  "
$ws.Range("D21").Formula = '="487.66"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -3.78%  "
$ws.Range("D22").Formula = '="9.38"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  -2.28%  "
$ws.Range("D23").Formula = '="0.723"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  -2.55%  "
$ws.Range("D24").Formula = '="85.44"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -0.74%  "
$ws.Range("E25").Value = "  -6.75%  "
$ws.Range("D26").Formula = '="0.0000138"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  -4.12%  "
$ws.Range("D27").Formula = '="12.16"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  -3.62%  "
$ws.Range("E28").Value = "  -3.65%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").Formula = '="2.93"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  -2.19%  "
$ws.Range("D31").Formula = '="2.35"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  -7.10%  "
$ws.Range("D32").Formula = '="31.47"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  -4.11%  "
$ws.Range("D33").Formula = '="7.62"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  -4.48%  "
$ws.Range("D34").Formula = '="3.846.13"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  -4.08%  "
$ws.Range("D35").Formula = '="0.107"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  -5.07%  "
$ws.Range("D36").Formula = '="3.652.40"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  -4.09%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").Formula = '="0.999"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  -4.80%  "
$ws.Range("D39").Formula = '="5.73"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  -6.33%  "
$ws.Range("E40").Value = "  -7.38%  "
$ws.Range("D41").Formula = '="0.321"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -4.55%  "
$ws.Range("D42").Formula = '="432.26"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  -8.85%  "
$ws.Range("D43").Formula = '="48.61"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -2.47%  "
$ws.Range("D44").Formula = '="1.93"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  -5.75%  "
$ws.Range("D45").Formula = '="2.78"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  -6.69%  "
$ws.Range("D46").Formula = '="8.39"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -1.74%  "
$ws.Range("D47").Formula = '="40.75"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -6.21%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").Formula = '="142.65"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  +1.65%  "
$ws.Range("D50").Formula = '="2.755.36"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  -5.94%  "
$ws.Range("D51").Formula = '="0.0347"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  -3.79%  "

$excel.CutCopyMode = $false
